# Update odds values in rows 3 and 4 of the sheet to match the refreshed
# FlashScore data snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("G3").Value = 1.48
$ws.Range("I3").Value = 7
$ws.Range("L3").Value = 7.5
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 7.5
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 2.5
$ws.Range("U3").Value = 2.5
$ws.Range("V3").Value = 1.5
$ws.Range("W3").Value = 5
$ws.Range("Z3").Value = 9.5
$ws.Range("AC3").Value = 7.5
$ws.Range("AE3").Value = 26
$ws.Range("AF3").Value = 101
$ws.Range("AJ3").Value = 23
$ws.Range("AL3").Value = 67
$ws.Range("AN3").Value = 3.2
$ws.Range("AQ3").Value = 23
$ws.Range("AS3").Value = 251
$ws.Range("AT3").Value = 2.5
$ws.Range("AU3").Value = 11
$ws.Range("AW3").Value = 8.5
$ws.Range("AY3").Value = 51
$ws.Range("AZ3").Value = 201
$ws.Range("BA3").Value = 251

# Row 4 updates
$ws.Range("J4").Value = 3.35
$ws.Range("K4").Value = 2.07
$ws.Range("L4").Value = 2.95
$ws.Range("O4").Value = 1.31
$ws.Range("P4").Value = 2.9
$ws.Range("Q4").Value = 1.95
$ws.Range("R4").Value = 1.75
$ws.Range("W4").Value = 8.25
$ws.Range("X4").Value = 14
$ws.Range("Y4").Value = 10.25
$ws.Range("AA4").Value = 25
$ws.Range("AB4").Value = 32
$ws.Range("AC4").Value = 8.75
$ws.Range("AH4").Value = 8.5
$ws.Range("AI4").Value = 13.5
$ws.Range("AJ4").Value = 9.25
$ws.Range("AK4").Value = 29
$ws.Range("AL4").Value = 20
$ws.Range("AM4").Value = 27
$ws.Range("AO4").Value = 15.5
$ws.Range("AP4").Value = 21
$ws.Range("AQ4").Value = 70
$ws.Range("AR4").Value = 100
$ws.Range("AS4").Value = 250
$ws.Range("AW4").Value = 4.5
$ws.Range("AX4").Value = 12.5
$ws.Range("AY4").Value = 17.5
$ws.Range("AZ4").Value = 50
$ws.Range("BA4").Value = 70
$ws.Range("BB4").Value = 175
